$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking price strings
# (e.g. "44.170.56", "1.90") are preserved exactly instead of being
# auto-converted into numbers by the input parser.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '44.170.56'
$ws.Range('E2').Value = '  -1.19%  '
$ws.Range('D3').Value = '2.248.60'
$ws.Range('E3').Value = '  -1.25%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '316.11'
$ws.Range('E5').Value = '  -1.60%  '
$ws.Range('D6').Value = '100.03'
$ws.Range('E6').Value = '  -6.02%  '
$ws.Range('D7').Value = '0.576'
$ws.Range('E7').Value = '  -3.18%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('D9').Value = '0.536'
$ws.Range('E9').Value = '  -6.69%  '
$ws.Range('D10').Value = '36.38'
$ws.Range('E10').Value = '  -6.12%  '
$ws.Range('D11').Value = '0.0826'
$ws.Range('E11').Value = '  -2.31%  '
$ws.Range('D12').Value = '7.39'
$ws.Range('E12').Value = '  -6.71%  '
$ws.Range('E13').Value = '  -2.86%  '
$ws.Range('D14').Value = '2.591.80'
$ws.Range('E14').Value = '  -1.20%  '
$ws.Range('D15').Value = '0.847'
$ws.Range('E15').Value = '  -4.37%  '
$ws.Range('D16').Value = '2.253.10'
$ws.Range('E16').Value = '  -1.09%  '
$ws.Range('D17').Value = '14.05'
$ws.Range('E17').Value = '  -4.13%  '
$ws.Range('D18').Value = '44.030.42'
$ws.Range('E18').Value = '  -1.06%  '
$ws.Range('E19').Value = '  -6.37%  '
$ws.Range('D20').Value = '0.0₃0984'
$ws.Range('E20').Value = '  -2.49%  '
$ws.Range('D21').Value = '6.34'
$ws.Range('E21').Value = '  -3.09%  '
$ws.Range('D22').Value = '65.75'
$ws.Range('E22').Value = '  -1.35%  '
$ws.Range('D23').Value = '238.35'
$ws.Range('E23').Value = '  -0.76%  '
$ws.Range('D24').Value = '2.99'
$ws.Range('E24').Value = '  -6.71%  '
$ws.Range('D25').Value = '2.04'
$ws.Range('E25').Value = '  -7.99%  '
$ws.Range('E26').Value = '  +0.33%  '
$ws.Range('D27').Value = '10.17'
$ws.Range('E27').Value = '  -0.46%  '
$ws.Range('E28').Value = '  -4.68%  '
$ws.Range('D29').Value = '36.73'
$ws.Range('E29').Value = '  -3.77%  '
$ws.Range('D30').Value = '6.01'
$ws.Range('E30').Value = '  -8.08%  '
$ws.Range('D31').Value = '20.10'
$ws.Range('E31').Value = '  -3.04%  '
$ws.Range('D32').Value = '155.96'
$ws.Range('E32').Value = '  -4.37%  '
$ws.Range('D33').Value = '0.0844'
$ws.Range('E33').Value = '  -5.16%  '
$ws.Range('D34').Value = '3.32'
$ws.Range('E34').Value = '  +5.49%  '
$ws.Range('E35').Value = '  -4.29%  '
$ws.Range('D36').Value = '1.90'
$ws.Range('E36').Value = '  -6.84%  '
$ws.Range('E37').Value = '  -8.04%  '
$ws.Range('E38').Value = '  -3.08%  '
$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').Value = '3.57'
$ws.Range('E39').Value = '  -9.33%  '
$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').Value = '15.37'
$ws.Range('E40').Value = '  -2.17%  '
$ws.Range('D41').Value = '3.98'
$ws.Range('E41').Value = '  -11.43%  '
$ws.Range('D42').Value = '0.0310'
$ws.Range('E42').Value = '  -6.05%  '
$ws.Range('D44').Value = '1.706.40'
$ws.Range('E44').Value = '  -4.01%  '
$ws.Range('D45').Value = '83.27'
$ws.Range('E45').Value = '  -4.94%  '
$ws.Range('E46').Value = '  -6.47%  '
$ws.Range('D47').Value = '5.21'
$ws.Range('E47').Value = '  -5.27%  '
$ws.Range('D48').Value = '102.10'
$ws.Range('E48').Value = '  -2.52%  '
$ws.Range('D49').Value = '71.51'
$ws.Range('E49').Value = '  -5.02%  '
$ws.Range('D50').Value = '56.83'
$ws.Range('E50').Value = '  -6.17%  '
$ws.Range('D51').Value = '1.62'
$ws.Range('E51').Value = '  -5.52%  '

# Restore the default cell style on column D (drop the temporary text
# number format) now that the literal text values are stored.
$ws.Range('D2:D51').Style = 'Normal'

Write-Output "done"
